$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3141.3333
$ws.Range("J32").Value = 4499.6665
$ws.Range("L32").Value = 4499.6665
$ws.Range("N32").Value = -5151.6665

$ws.Range("H38").Value = 1472.5
$ws.Range("I38").Value = 1472.5
$ws.Range("K38").Value = 4417.5
$ws.Range("M38").Value = -4045.5

$ws.Range("H80").Value = 754.8182
$ws.Range("J80").Value = 826.875
$ws.Range("L80").Value = 2480.625
$ws.Range("N80").Value = -4476.625

$ws.Range("H83").Value = 754.8182
$ws.Range("J83").Value = 826.875
$ws.Range("L83").Value = 7441.875
$ws.Range("N83").Value = -17425.875

$ws.Range("H86").Value = 2004043.8
$ws.Range("I86").Value = 5002109.5
$ws.Range("K86").Value = 5002109.5
$ws.Range("M86").Value = -5000986.5

$ws.Range("H89").Value = 2004043.8
$ws.Range("I89").Value = 5002109.5
$ws.Range("K89").Value = 25010547.5
$ws.Range("M89").Value = -25004931.5

$ws.Range("H112").Value = 2130.2307
$ws.Range("J112").Value = 2130.2307
$ws.Range("L112").Value = 6390.6921
$ws.Range("N112").Value = -8606.6921

$ws.Range("H116").Value = 8340.299999999999
$ws.Range("I116").Value = 8057.5713
$ws.Range("J116").Value = 9000
$ws.Range("K116").Value = 8057.5713
$ws.Range("L116").Value = 9000
$ws.Range("M116").Value = -4615.5713
$ws.Range("N116").Value = -15884

$ws.Range("H137").Value = 3999.2307
$ws.Range("I137").Value = 2608.0356
$ws.Range("K137").Value = 7824.1068
$ws.Range("M137").Value = -5274.1068


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 26009208
$ws.Range("I74").Value = 50004680
$ws.Range("K74").Value = 50004680
$ws.Range("M74").Value = -50003806

$ws.Range("H77").Value = 26009208
$ws.Range("I77").Value = 50004680
$ws.Range("K77").Value = 250023400
$ws.Range("M77").Value = -250019032

$ws.Range("H122").Value = 3249.5
$ws.Range("I122").Value = 2571.25
$ws.Range("K122").Value = 7713.75
$ws.Range("M122").Value = -5263.75

$ws.Range("H126").Value = 7400
$ws.Range("I126").Value = 7400
$ws.Range("K126").Value = 22200
$ws.Range("M126").Value = -19730

$ws.Range("H132").Value = 8845.76
$ws.Range("I132").Value = 6287
$ws.Range("K132").Value = 18861
$ws.Range("M132").Value = -16331


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3732.0417
$ws.Range("I99").Value = 3084.9
$ws.Range("J99").Value = 4194.2856
$ws.Range("K99").Value = 3084.9
$ws.Range("L99").Value = 4194.2856
$ws.Range("M99").Value = -1586.9
$ws.Range("N99").Value = -7190.2856

$ws.Range("H105").Value = 2670.1052
$ws.Range("I105").Value = 3275
$ws.Range("K105").Value = 3275
$ws.Range("M105").Value = -1528

$ws.Range("H107").Value = 1640.64
$ws.Range("I107").Value = 1707.6818
$ws.Range("K107").Value = 1707.6818
$ws.Range("M107").Value = 212.3181999999999

$ws.Range("H134").Value = 1001254.9
$ws.Range("I134").Value = 1066.875
$ws.Range("K134").Value = 3200.625
$ws.Range("M134").Value = -665.625


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 576872.3
$ws.Range("I31").Value = 8852.087
$ws.Range("K31").Value = 8852.087
$ws.Range("M31").Value = -8557.087

$ws.Range("H34").Value = 576872.3
$ws.Range("I34").Value = 8852.087
$ws.Range("K34").Value = 8852.087
$ws.Range("M34").Value = -8650.087

$ws.Range("H58").Value = 2535.6667
$ws.Range("I58").Value = 2129.818
$ws.Range("K58").Value = 2129.818
$ws.Range("M58").Value = -1926.818

$ws.Range("H86").Value = 9952
$ws.Range("I86").Value = 9952
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 9952
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -8829
$ws.Range("N86").ClearContents() | Out-Null

$ws.Range("H89").Value = 9952
$ws.Range("I89").Value = 9952
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 49760
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -44144
$ws.Range("N89").ClearContents() | Out-Null

$ws.Range("H99").Value = 3999.5
$ws.Range("I99").Value = 3999.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3999.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2501.5
$ws.Range("N99").ClearContents() | Out-Null

$ws.Range("H108").Value = 90975.60000000001
$ws.Range("J108").Value = 90975.60000000001
$ws.Range("L108").Value = 90975.60000000001
$ws.Range("N108").Value = -98655.60000000001

$ws.Range("H126").Value = 3999.5
$ws.Range("I126").Value = 3999.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11998.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -9528.5
$ws.Range("N126").ClearContents() | Out-Null

$ws.Range("H132").Value = 2055.8125
$ws.Range("I132").Value = 1926.2667
$ws.Range("K132").Value = 5778.800099999999
$ws.Range("M132").Value = -3248.800099999999

$ws.Range("H134").Value = 3891.5264
$ws.Range("I134").Value = 2947.8667
$ws.Range("J134").Value = 7430.25
$ws.Range("K134").Value = 8843.6001
$ws.Range("L134").Value = 22290.75
$ws.Range("M134").Value = -6308.6001
$ws.Range("N134").Value = -27360.75

$ws.Range("H136").Value = 2535.6667
$ws.Range("I136").Value = 2129.818
$ws.Range("K136").Value = 6389.454000000001
$ws.Range("M136").Value = -3839.454000000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1309.5555
$ws.Range("I46").Value = 826.5714
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 2479.7142
$ws.Range("L46").Value = 9000
$ws.Range("M46").Value = -2388.7142
$ws.Range("N46").Value = -9182

$ws.Range("H80").Value = 4349.231
$ws.Range("I80").Value = 4402
$ws.Range("J80").Value = 4344.8335
$ws.Range("K80").Value = 13206
$ws.Range("L80").Value = 13034.5005
$ws.Range("M80").Value = -12270
$ws.Range("N80").Value = -14906.5005

$ws.Range("H83").Value = 4349.231
$ws.Range("I83").Value = 4402
$ws.Range("J83").Value = 4344.8335
$ws.Range("K83").Value = 39618
$ws.Range("L83").Value = 39103.5015
$ws.Range("M83").Value = -34938
$ws.Range("N83").Value = -48463.5015


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7435.5386
$ws.Range("I70").Value = 6857.875
$ws.Range("K70").Value = 6857.875
$ws.Range("M70").Value = -6587.875

$ws.Range("H73").Value = 7435.5386
$ws.Range("I73").Value = 6857.875
$ws.Range("K73").Value = 6857.875
$ws.Range("M73").Value = -5921.875

$ws.Range("H92").Value = 25050
$ws.Range("J92").Value = 25050
$ws.Range("L92").Value = 25050
$ws.Range("N92").Value = -28794

$ws.Range("H107").Value = 377.3889
$ws.Range("I107").Value = 313.35715
$ws.Range("K107").Value = 313.35715
$ws.Range("M107").Value = 1606.64285

$ws.Range("H122").Value = 2293.8823
$ws.Range("I122").Value = 2066.4
$ws.Range("K122").Value = 6199.200000000001
$ws.Range("M122").Value = -3749.200000000001

$ws.Range("H126").Value = 3490.1875
$ws.Range("I126").Value = 3072.75
$ws.Range("J126").Value = 3629.3333
$ws.Range("K126").Value = 9218.25
$ws.Range("L126").Value = 10887.9999
$ws.Range("M126").Value = -6748.25
$ws.Range("N126").Value = -15827.9999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31106
$ws.Range("I7").Value = 3086.8333
$ws.Range("J7").Value = 82833.69500000001
$ws.Range("K7").Value = 3086.8333
$ws.Range("L7").Value = 82833.69500000001
$ws.Range("M7").Value = -2974.8333
$ws.Range("N7").Value = -83057.69500000001

$ws.Range("H40").Value = 3884.2
$ws.Range("I40").Value = 3216.0527
$ws.Range("K40").Value = 3216.0527
$ws.Range("M40").Value = -3080.0527

$ws.Range("H42").Value = 30000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -31126

$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30294

$ws.Range("H61").Value = 1604.4445
$ws.Range("I61").Value = 1348.5714
$ws.Range("K61").Value = 1348.5714
$ws.Range("M61").Value = -1146.5714

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents() | Out-Null
$ws.Range("N98").ClearContents() | Out-Null

$ws.Range("H113").Value = 1604.4445
$ws.Range("I113").Value = 1348.5714
$ws.Range("K113").Value = 1348.5714
$ws.Range("M113").Value = 821.4286

$ws.Range("H122").Value = 5046.2334
$ws.Range("I122").Value = 4676.591
$ws.Range("K122").Value = 14029.773
$ws.Range("M122").Value = -11579.773

$ws.Range("H126").Value = 31106
$ws.Range("I126").Value = 3086.8333
$ws.Range("J126").Value = 82833.69500000001
$ws.Range("K126").Value = 9260.499899999999
$ws.Range("L126").Value = 248501.085
$ws.Range("M126").Value = -6790.499899999999
$ws.Range("N126").Value = -253441.085

$ws.Range("H132").Value = 779955.9399999999
$ws.Range("I132").Value = 12642.2
$ws.Range("K132").Value = 37926.60000000001
$ws.Range("M132").Value = -35396.60000000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2375.5
$ws.Range("I126").Value = 1501.3334
$ws.Range("J126").Value = 4998
$ws.Range("K126").Value = 4504.0002
$ws.Range("L126").Value = 14994
$ws.Range("M126").Value = -2034.0002
$ws.Range("N126").Value = -19934

$ws.Range("H136").Value = 4275.3687
$ws.Range("I136").Value = 4597.3125
$ws.Range("J136").Value = 2558.3333
$ws.Range("K136").Value = 13791.9375
$ws.Range("L136").Value = 7674.999899999999
$ws.Range("M136").Value = -11241.9375
$ws.Range("N136").Value = -12774.9999

